$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.762.09"
$ws.Range("E2").Value = '  +7.29%  '
$ws.Range("D3").Value = "'1.737.99"
$ws.Range("E3").Value = '  +3.53%  '
$ws.Range("D4").Value = "'0.9981"
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("D5").Value = "'333.54"
$ws.Range("E5").Value = '  +1.57%  '
$ws.Range("D6").Value = "'0.9965"
$ws.Range("E6").Value = '  -0.16%  '
$ws.Range("D7").Value = "'0.3738"
$ws.Range("E7").Value = '  +2.34%  '
$ws.Range("D8").Value = "'0.3395"
$ws.Range("E8").Value = '  +4.13%  '
$ws.Range("D9").Value = "'48.22"
$ws.Range("E9").Value = '  +2.38%  '
$ws.Range("D10").Value = "'1.187"
$ws.Range("E10").Value = '  +3.40%  '
$ws.Range("D11").Value = "'0.07455"
$ws.Range("E11").Value = '  +5.09%  '
$ws.Range("D12").Value = "'0.9968"
$ws.Range("E12").Value = '  -0.14%  '
$ws.Range("D13").Value = "'6.418"
$ws.Range("E13").Value = '  +5.27%  '
$ws.Range("E14").Value = '  +3.34%  '
$ws.Range("E15").Value = '  +6.66%  '
$ws.Range("D16").Value = "'1.734.44"
$ws.Range("E16").Value = '  +3.43%  '
$ws.Range("D17").Value = "'0.00001077"
$ws.Range("E17").Value = '  +2.03%  '
$ws.Range("D18").Value = "'0.06729"
$ws.Range("E18").Value = '  +1.86%  '
$ws.Range("D19").Value = "'82.65"
$ws.Range("E19").Value = '  +4.41%  '
$ws.Range("D20").Value = "'0.9962"
$ws.Range("E20").Value = '  -0.12%  '
$ws.Range("D21").Value = "'16.68"
$ws.Range("E21").Value = '  +4.36%  '
$ws.Range("D22").Value = "'6.221"
$ws.Range("E22").Value = '  +4.53%  '
$ws.Range("D23").Value = "'12.77"
$ws.Range("E23").Value = '  +2.02%  '
$ws.Range("D24").Value = "'26.733.97"
$ws.Range("E24").Value = '  +7.22%  '
$ws.Range("D25").Value = "'2.443"
$ws.Range("E25").Value = '  -0.39%  '
$ws.Range("D26").Value = "'1.487"
$ws.Range("E26").Value = '  +24.88%  '
$ws.Range("D27").Value = "'2.425"
$ws.Range("E27").Value = '  -0.07%  '
$ws.Range("D28").Value = "'151.42"
$ws.Range("E28").Value = '  +1.88%  '
$ws.Range("D29").Value = "'19.55"
$ws.Range("E29").Value = '  +3.81%  '
$ws.Range("D30").Value = "'1.930.47"
$ws.Range("E30").Value = '  +3.69%  '
$ws.Range("D31").Value = "'132.31"
$ws.Range("E31").Value = '  +4.90%  '
$ws.Range("D32").Value = "'4.096"
$ws.Range("E32").Value = '  +0.73%  '
$ws.Range("D33").Value = "'6.049"
$ws.Range("E33").Value = '  +4.57%  '
$ws.Range("D34").Value = "'0.08639"
$ws.Range("E34").Value = '  +2.06%  '
$ws.Range("D35").Value = "'1.700"
$ws.Range("E35").Value = '  +3.14%  '
$ws.Range("E36").Value = '  +4.68%  '
$ws.Range("D37").Value = "'5.399"
$ws.Range("E37").Value = '  +4.12%  '
$ws.Range("D38").Value = "'0.02350"
$ws.Range("E38").Value = '  +3.35%  '
$ws.Range("D39").Value = "'0.2177"
$ws.Range("E39").Value = '  +4.03%  '
$ws.Range("D40").Value = "'0.06235"
$ws.Range("E40").Value = '  +2.56%  '
$ws.Range("D41").Value = "'8.461"
$ws.Range("E41").Value = '  +2.00%  '
$ws.Range("D42").Value = "'1.226"
$ws.Range("E42").Value = '  -0.11%  '
$ws.Range("D43").Value = "'0.6261"
$ws.Range("E43").Value = '  +4.72%  '
$ws.Range("D44").Value = "'14.29"
$ws.Range("E44").Value = '  +5.08%  '
$ws.Range("D45").Value = "'0.9953"
$ws.Range("E45").Value = '  -0.10%  '
$ws.Range("E46").Value = '  +2.25%  '
$ws.Range("D47").Value = "'0.6098"
$ws.Range("E47").Value = '  +6.40%  '
$ws.Range("E48").Value = '  +2.76%  '
$ws.Range("D49").Value = "'2.062"
$ws.Range("E49").Value = '  +4.88%  '
$ws.Range("D50").Value = "'0.07212"
$ws.Range("E50").Value = '  +2.69%  '
$ws.Range("D51").Value = "'77.64"
$ws.Range("E51").Value = '  +3.38%  '
